# Fruta / hortaliza, semanal
# Insert 2 new weekly report rows (Doctor Davis variety, Primera & Segunda
# quality) right above the existing data block that starts at row 281. This
# pushes all the previously-recorded rows down by two positions and grows
# the used range from A1:T348 to A1:T350.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 281 - everything that used
# to live at 281..348 now lives at 283..350.
$ws.Rows("281:282").Insert()

# ---- Row 281: Doctor Davis / Primera -------------------------------------
$ws.Cells.Item(281, 1).Value  = 11
$ws.Cells.Item(281, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(281, 3).Value  = "Bíobío"
$ws.Cells.Item(281, 4).Value  = 44985
$ws.Cells.Item(281, 5).Value  = 8
$ws.Cells.Item(281, 6).Value  = "Fruta"
$ws.Cells.Item(281, 7).Value  = 100103
$ws.Cells.Item(281, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(281, 9).Value  = 100103004
$ws.Cells.Item(281, 10).Value = "Durazno"
$ws.Cells.Item(281, 11).Value = "Doctor Davis"
$ws.Cells.Item(281, 12).Value = "Primera"
$ws.Cells.Item(281, 13).Value = 150
$ws.Cells.Item(281, 14).Value = 14000
$ws.Cells.Item(281, 15).Value = 14000
$ws.Cells.Item(281, 16).Value = 14000
$ws.Cells.Item(281, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(281, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(281, 19).Value = 875
$ws.Cells.Item(281, 20).Value = 16

# ---- Row 282: Doctor Davis / Segunda --------------------------------------
$ws.Cells.Item(282, 1).Value  = 11
$ws.Cells.Item(282, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(282, 3).Value  = "Bíobío"
$ws.Cells.Item(282, 4).Value  = 44985
$ws.Cells.Item(282, 5).Value  = 8
$ws.Cells.Item(282, 6).Value  = "Fruta"
$ws.Cells.Item(282, 7).Value  = 100103
$ws.Cells.Item(282, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(282, 9).Value  = 100103004
$ws.Cells.Item(282, 10).Value = "Durazno"
$ws.Cells.Item(282, 11).Value = "Doctor Davis"
$ws.Cells.Item(282, 12).Value = "Segunda"
$ws.Cells.Item(282, 13).Value = 100
$ws.Cells.Item(282, 14).Value = 12000
$ws.Cells.Item(282, 15).Value = 12000
$ws.Cells.Item(282, 16).Value = 12000
$ws.Cells.Item(282, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(282, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(282, 19).Value = 750
$ws.Cells.Item(282, 20).Value = 16

# Make sure the date column keeps the existing date number-format that the
# surrounding rows use (Insert() already propagates formatting from the row
# above, but re-assert the number format explicitly to be safe).
$dateFmt = $ws.Range("D283").NumberFormat
$ws.Range("D281").NumberFormat = $dateFmt
$ws.Range("D282").NumberFormat = $dateFmt
